$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = -13
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = -2
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = -2
